$d = $word.ActiveDocument

$replacements = @(
    @("779÷8=", "641÷3="),
    @("520÷5=", "464÷4="),
    @("573÷6=", "723÷6="),
    @("402÷5=", "285÷6="),
    @("400÷9=", "612÷5="),
    @("762÷9=", "417÷6="),
    @("843÷9=", "913÷9="),
    @("479÷5=", "814÷3="),
    @("304÷9=", "143÷7="),
    @("699÷3=", "810÷6="),
    @("871÷5=", "466÷4="),
    @("673÷5=", "931÷2="),
    @("847÷5=", "178÷8="),
    @("374÷2=", "967÷4="),
    @("605÷9=", "212÷9="),
    @("756÷8=", "600÷2="),
    @("260÷4=", "509÷2="),
    @("588÷8=", "900÷2="),
    @("179÷8=", "696÷7="),
    @("576÷2=", "225÷3="),
    @("749÷9=", "954÷6="),
    @("586÷7=", "676÷5="),
    @("564÷2=", "561÷2="),
    @("132÷8=", "496÷4="),
    @("453÷7=", "834÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
